# Refresh the live crypto price/volume snapshot on Sheet1.
# Each row holds one coin: B=Name, C=Link, D=Price, E=Volume(1h)).
# Price cells that look like plain numbers are written with a leading
# apostrophe so Excel keeps them as text (preserving trailing zeros
# and the "thousand.thousand.decimal" display format already used in
# this sheet) instead of silently re-parsing them as floats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '70.920.44'
$ws.Range('E2').Value = '  +0.51%  '

# Row 3
$ws.Range('D3').Value = '3.539.58'
$ws.Range('E3').Value = '  -0.66%  '

# Row 4
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.10%  '

# Row 5
$ws.Range('D5').Value = '''625.49'
$ws.Range('E5').Value = '  +2.43%  '

# Row 6
$ws.Range('D6').Value = '''175.09'
$ws.Range('E6').Value = '  +0.96%  '

# Row 7
$ws.Range('D7').Value = '3.537.72'
$ws.Range('E7').Value = '  -0.63%  '

# Row 8
$ws.Range('E8').Value = '  -1.11%  '

# Row 9
$ws.Range('D9').Value = '''0.999'
$ws.Range('E9').Value = '  -0.05%  '

# Row 10
$ws.Range('E10').Value = '  +1.26%  '

# Row 11
$ws.Range('D11').Value = '''7.20'
$ws.Range('E11').Value = '  -5.74%  '

# Row 12
$ws.Range('D12').Value = '''0.588'
$ws.Range('E12').Value = '  +0.19%  '

# Row 13
$ws.Range('D13').Value = '''46.84'
$ws.Range('E13').Value = '  +0.34%  '

# Row 14
$ws.Range('E14').Value = '  +0.06%  '

# Row 15
$ws.Range('D15').Value = '4.109.00'
$ws.Range('E15').Value = '  -0.90%  '

# Row 16
$ws.Range('D16').Value = '''8.45'
$ws.Range('E16').Value = '  +0.70%  '

# Row 17
$ws.Range('D17').Value = '''610.47'
$ws.Range('E17').Value = '  -0.41%  '

# Row 18
$ws.Range('D18').Value = '3.541.05'
$ws.Range('E18').Value = '  -0.99%  '

# Row 19
$ws.Range('D19').Value = '70.971.62'
$ws.Range('E19').Value = '  +0.47%  '

# Row 20
$ws.Range('E20').Value = '  +1.29%  '

# Row 21
$ws.Range('D21').Value = '''17.84'
$ws.Range('E21').Value = '  +2.47%  '

# Row 22
$ws.Range('D22').Value = '''0.889'
$ws.Range('E22').Value = '  +0.35%  '

# Row 23
$ws.Range('D23').Value = '''9.08'
$ws.Range('E23').Value = '  -3.71%  '

# Row 24
$ws.Range('D24').Value = '''15.76'
$ws.Range('E24').Value = '  -2.07%  '

# Row 25
$ws.Range('D25').Value = '''98.44'
$ws.Range('E25').Value = '  +1.35%  '

# Row 26
$ws.Range('D26').Value = '''3.81'
$ws.Range('E26').Value = '  -0.42%  '

# Row 28
$ws.Range('D28').Value = '''2.60'
$ws.Range('E28').Value = '  -1.08%  '

# Row 29
$ws.Range('D29').Value = '''34.03'
$ws.Range('E29').Value = '  +1.61%  '

# Row 30
$ws.Range('D30').Value = '''9.17'
$ws.Range('E30').Value = '  +0.88%  '

# Row 31
$ws.Range('D31').Value = '''3.07'
$ws.Range('E31').Value = '  +0.51%  '

# Row 32
$ws.Range('D32').Value = '''8.19'
$ws.Range('E32').Value = '  -3.76%  '

# Row 33
$ws.Range('E33').Value = '  +0.59%  '

# Row 34
$ws.Range('D34').Value = '''6.90'
$ws.Range('E34').Value = '  -1.26%  '

# Row 35
$ws.Range('D35').Value = '''632.65'
$ws.Range('E35').Value = '  +9.67%  '

# Row 36
$ws.Range('E36').Value = '  -0.96%  '

# Row 37
$ws.Range('D37').Value = '''10.88'
$ws.Range('E37').Value = '  +0.47%  '

# Row 38
$ws.Range('D38').Value = '''3.52'
$ws.Range('E38').Value = '  -3.85%  '

# Row 39
$ws.Range('D39').Value = '''0.0477'
$ws.Range('E39').Value = '  -1.55%  '

# Row 40
$ws.Range('D40').Value = '''57.02'
$ws.Range('E40').Value = '  -0.66%  '

# Row 41
$ws.Range('E41').Value = '  +0.25%  '

# Row 42
$ws.Range('E42').Value = '  +2.17%  '

# Row 43
$ws.Range('D43').Value = '0.0₃0743'
$ws.Range('E43').Value = '  +5.30%  '

# Row 44
$ws.Range('D44').Value = '3.368.02'
$ws.Range('E44').Value = '  -0.52%  '

# Row 45
$ws.Range('E45').Value = '  +0.39%  '

# Row 46
$ws.Range('E46').Value = '  -1.94%  '

# Row 47
$ws.Range('D47').Value = '''32.31'
$ws.Range('E47').Value = '  -2.93%  '

# Row 48
$ws.Range('D48').Value = '''2.58'
$ws.Range('E48').Value = '  -1.12%  '

# Row 49
$ws.Range('E49').Value = '  +0.49%  '

# Row 50
$ws.Range('D50').Value = '''133.14'
$ws.Range('E50').Value = '  -0.61%  '

# Row 51
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.156'
$ws.Range('E51').Value = '  +5.08%  '
